# Swap the order of the first two comma-separated names in the
# "Recorded By" column (column G) wherever the entry begins with
# "System, ...". E.g. "System, foo@bar.com" -> "foo@bar.com, System"
# and "System, foo@bar.com, system" -> "foo@bar.com, System, system".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value
    $parts = $text -split ', '

    if ($parts.Length -ge 2 -and $parts[0] -eq 'System') {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $cell.Value2 = [string]::Join(', ', $parts)
    }
}
